$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to the changed Price/Volume cells so Excel
# keeps numeric-looking strings (e.g. "1.001") stored as text, matching
# the original inline-string cell typing instead of auto-converting to numbers.
$data = [ordered]@{
    'D2' = '30.663.06'
    'E2' = '  +0.54%  '
    'D3' = '1.877.03'
    'E3' = '  -0.93%  '
    'D4' = '1.001'
    'E4' = '  +0.02%  '
    'D5' = '239.08'
    'E5' = '  -0.11%  '
    'E6' = '  +0.03%  '
    'D7' = '0.4797'
    'E7' = '  -1.15%  '
    'D8' = '0.2823'
    'E8' = '  -2.65%  '
    'D9' = '0.06499'
    'E9' = '  -1.83%  '
    'D10' = '1.882.81'
    'E10' = '  -0.46%  '
    'D11' = '0.07466'
    'E11' = '  +0.63%  '
    'D12' = '16.48'
    'E12' = '  -2.41%  '
    'D13' = '5.085'
    'E13' = '  -2.32%  '
    'D14' = '87.85'
    'E14' = '  -1.22%  '
    'D15' = '0.6629'
    'E15' = '  -0.10%  '
    'D16' = '30.604.75'
    'E16' = '  +0.50%  '
    'D17' = '13.25'
    'E17' = '  -2.29%  '
    'D18' = '1.002'
    'E18' = '  -0.08%  '
    'B19' = 'WrappedliquidstakedEther2.0'
    'C19' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D19' = '2.203.29'
    'E19' = '  +3.03%  '
    'B20' = 'ShibaInu'
    'C20' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D20' = '0.000007564'
    'E20' = '  -2.99%  '
    'D21' = '226.72'
    'E21' = '  +1.64%  '
    'D22' = '1.002'
    'E22' = '  +0.02%  '
    'D23' = '5.272'
    'E23' = '  -2.24%  '
    'D24' = '6.142'
    'E24' = '  -1.64%  '
    'B25' = 'Cosmos'
    'C25' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D25' = '9.291'
    'E25' = '  -1.25%  '
    'B26' = 'Monero'
    'C26' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D26' = '167.85'
    'E26' = '  +2.65%  '
    'D27' = '18.52'
    'D28' = '1.929'
    'E28' = '  -1.21%  '
    'D29' = '1.406'
    'E29' = '  -2.97%  '
    'D30' = '0.09684'
    'E30' = '  +4.76%  '
    'D31' = '4.333'
    'E31' = '  -0.23%  '
    'D32' = '3.998'
    'E32' = '  -1.29%  '
    'D33' = '0.05053'
    'E33' = '  -0.57%  '
    'D34' = '1.214'
    'E34' = '  +4.42%  '
    'D35' = '0.7457'
    'E35' = '  -2.41%  '
    'D36' = '2.716'
    'E36' = '  +0.63%  '
    'D37' = '0.01859'
    'E37' = '  -1.22%  '
    'D38' = '2.637'
    'E38' = '  -0.33%  '
    'D39' = '0.9129'
    'E39' = '  -0.79%  '
    'D40' = '2.069'
    'E40' = '  -1.16%  '
    'D41' = '105.99'
    'E41' = '  -0.59%  '
    'D42' = '0.4264'
    'E42' = '  -2.44%  '
    'D43' = '5.771'
    'E43' = '  -3.20%  '
    'D44' = '0.9990'
    'E44' = '  -0.44%  '
    'D45' = '7.326'
    'E45' = '  -4.20%  '
    'D46' = '0.1290'
    'E46' = '  -3.15%  '
    'D47' = '64.09'
    'E47' = '  -2.51%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '8.916'
    'E48' = '  -0.62%  '
    'B49' = 'NEARProtocol'
    'C49' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D49' = '1.467'
    'E49' = '  -8.64%  '
    'D50' = '33.69'
    'E50' = '  -2.79%  '
    'D51' = '0.05655'
    'E51' = '  -1.16%  '
}

foreach ($addr in $data.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $data[$addr]
}
